$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Step A: trim paragraph 31 ("31 Тепер якщо вам непотрібна програма Skype...").
#   - extend "...способом " run text to merge with the text that used to sit
#     after the _GoBack bookmark, ending with "і до цього." (no trailing
#     spaces any more).
#   - this also removes the (now redundant) _GoBack bookmark and the
#     trailing whitespace-only runs, because they all sit inside the
#     replaced span.
# -----------------------------------------------------------------------
$p31 = $d.Paragraphs.Item(31)
$r31 = $p31.Range
$oldTail = "видаляємо її тим самим способом який ми використовували і до цього.         "
$newTail = "видаляємо її тим самим способом який ми використовували і до цього."
$r31.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null

# -----------------------------------------------------------------------
# Step B: old paragraph 35 ("35. Сьогодні давайте...") -> renumber to "37".
#   Replacing the numeral run's text also drops the lastRenderedPageBreak
#   marker that used to sit in front of it (it is a pure render artefact).
# -----------------------------------------------------------------------
$p35 = $d.Paragraphs.Item(35)
$num35 = $d.Range($p35.Range.Start, $p35.Range.Start + 2)
$num35.Find.Execute("35", $true, $false, $false, $false, $false, $true, 1, $false, "37", 2) | Out-Null

# -----------------------------------------------------------------------
# Step C: old paragraph 34 ("34. Тепер давайте зайдемо...") -> renumber "36".
# -----------------------------------------------------------------------
$p34 = $d.Paragraphs.Item(34)
$num34 = $d.Range($p34.Range.Start, $p34.Range.Start + 2)
$num34.Find.Execute("34", $true, $false, $false, $false, $false, $true, 1, $false, "36", 2) | Out-Null

# -----------------------------------------------------------------------
# Step D: old paragraph 33 ("33. Слідом дану дію...") -> renumber "34", then
#   insert a brand new paragraph "35." right after it (with the
#   lastRenderedPageBreak + _GoBack bookmark moved here, as in the diff).
# -----------------------------------------------------------------------
$p33 = $d.Paragraphs.Item(33)
$num33 = $d.Range($p33.Range.Start, $p33.Range.Start + 2)
$num33.Find.Execute("33", $true, $false, $false, $false, $false, $true, 1, $false, "34", 2) | Out-Null

$p33now = $d.Paragraphs.Item(33)
$p33now.Range.InsertParagraphAfter()
$newPara35 = $d.Paragraphs.Item(34)
$xml35 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:spacing w:before=`"240`"/></w:pPr><w:r><w:rPr><w:lang w:val=`"ru-RU`"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">35. </w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:r><w:t xml:space=`"preserve`">За цим давайте знову спустимось в низ нашого меню вже знайомим нам способом.  </w:t></w:r><w:bookmarkEnd w:id=`"0`"/></w:p>"
$newPara35.Range.InsertXML($xml35) | Out-Null

# -----------------------------------------------------------------------
# Step E: old paragraph 32 ("32. Далі заходимо в папку система Windows...")
#   -> renumber "32" to "33" (only the second digit run changes: 2 -> 3).
# -----------------------------------------------------------------------
$p32 = $d.Paragraphs.Item(32)
$num32 = $d.Range($p32.Range.Start + 1, $p32.Range.Start + 2)
$num32.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "3", 2) | Out-Null

# -----------------------------------------------------------------------
# Step F: insert the brand new paragraph "32" (about scrolling down the
#   start menu) right after paragraph 31.
# -----------------------------------------------------------------------
$p31now = $d.Paragraphs.Item(31)
$p31now.Range.InsertParagraphAfter()
$newPara32 = $d.Paragraphs.Item(32)
$xml32 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:spacing w:before=`"240`"/></w:pPr><w:r><w:rPr><w:lang w:val=`"ru-RU`"/></w:rPr><w:t xml:space=`"preserve`">32 </w:t></w:r><w:r><w:t xml:space=`"preserve`">За цим давайте спустимось в низ нашого меню, для цього прокручуємо коліщатко нашої миші на себе.          </w:t></w:r><w:r><w:rPr><w:lang w:val=`"ru-RU`"/></w:rPr><w:t xml:space=`"preserve`">  </w:t></w:r><w:r><w:t xml:space=`"preserve`">  </w:t></w:r></w:p>"
$newPara32.Range.InsertXML($xml32) | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
